# Auto-generated edit script applying the diff changes to Kujata_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1701.5
$ws.Range("I103").Value = 504
$ws.Range("J103").Value = 1941
$ws.Range("K103").Value = 1512
$ws.Range("L103").Value = 5823
$ws.Range("M103").Value = -926
$ws.Range("N103").Value = -6995
$ws.Range("H125").Value = 1776.1428
$ws.Range("I125").Value = 1730.5
$ws.Range("J125").Value = 1794.4
$ws.Range("K125").Value = 15574.5
$ws.Range("L125").Value = 16149.6
$ws.Range("M125").Value = -13114.5
$ws.Range("N125").Value = -21069.6
$ws.Range("H135").Value = 31250646
$ws.Range("I135").Value = 290.68
$ws.Range("K135").Value = 2616.12
$ws.Range("M135").Value = -81.11999999999989
$ws.Range("H137").Value = 1144.7073
$ws.Range("I137").Value = 791.8261
$ws.Range("K137").Value = 2375.4783
$ws.Range("M137").Value = 174.5217000000002
$ws.Range("H141").Value = 1004.0909
$ws.Range("I141").Value = 759.44446
$ws.Range("J141").Value = 2105
$ws.Range("K141").Value = 2278.33338
$ws.Range("L141").Value = 6315
$ws.Range("M141").Value = 2901.66662
$ws.Range("N141").Value = -16675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17391
$ws.Range("I2").Value = 633.25
$ws.Range("J2").Value = 50906.5
$ws.Range("K2").Value = 633.25
$ws.Range("L2").Value = 50906.5
$ws.Range("M2").Value = -520.25
$ws.Range("N2").Value = -51132.5
$ws.Range("H32").Value = 3567.2957
$ws.Range("I32").Value = 3276.918
$ws.Range("K32").Value = 3276.918
$ws.Range("M32").Value = -2989.918
$ws.Range("H61").Value = 1542.5769
$ws.Range("I61").Value = 1372.7142
$ws.Range("J61").Value = 2256
$ws.Range("K61").Value = 1372.7142
$ws.Range("L61").Value = 2256
$ws.Range("M61").Value = -1160.7142
$ws.Range("N61").Value = -2680
$ws.Range("H74").Value = 1044.75
$ws.Range("I74").Value = 762.3158
$ws.Range("J74").Value = 2118
$ws.Range("K74").Value = 762.3158
$ws.Range("L74").Value = 2118
$ws.Range("M74").Value = 111.6842
$ws.Range("N74").Value = -3866
$ws.Range("H77").Value = 1044.75
$ws.Range("I77").Value = 762.3158
$ws.Range("J77").Value = 2118
$ws.Range("K77").Value = 3811.579
$ws.Range("L77").Value = 10590
$ws.Range("M77").Value = 556.4210000000003
$ws.Range("N77").Value = -19326
$ws.Range("H88").Value = 2986.0715
$ws.Range("I88").Value = 2452.5
$ws.Range("J88").Value = 3075
$ws.Range("K88").Value = 2452.5
$ws.Range("L88").Value = 3075
$ws.Range("M88").Value = -2046.5
$ws.Range("N88").Value = -3887
$ws.Range("H91").Value = 2986.0715
$ws.Range("I91").Value = 2452.5
$ws.Range("J91").Value = 3075
$ws.Range("K91").Value = 2452.5
$ws.Range("L91").Value = 3075
$ws.Range("M91").Value = -1048.5
$ws.Range("N91").Value = -5883
$ws.Range("H97").Value = 406.9091
$ws.Range("I97").Value = 406.9091
$ws.Range("K97").Value = 406.9091
$ws.Range("M97").Value = 89.09089999999998
$ws.Range("H102").Value = 166666670
$ws.Range("I102").Value = 166666670
$ws.Range("K102").Value = 166666670
$ws.Range("M102").Value = -166665048
$ws.Range("H116").Value = 17391
$ws.Range("I116").Value = 633.25
$ws.Range("J116").Value = 50906.5
$ws.Range("K116").Value = 633.25
$ws.Range("L116").Value = 50906.5
$ws.Range("M116").Value = 1660.75
$ws.Range("N116").Value = -55494.5
$ws.Range("H136").Value = 1542.5769
$ws.Range("I136").Value = 1372.7142
$ws.Range("J136").Value = 2256
$ws.Range("K136").Value = 4118.142599999999
$ws.Range("L136").Value = 6768
$ws.Range("M136").Value = -1568.142599999999
$ws.Range("N136").Value = -11868

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17391
$ws.Range("I3").Value = 633.25
$ws.Range("J3").Value = 50906.5
$ws.Range("K3").Value = 633.25
$ws.Range("L3").Value = 50906.5
$ws.Range("M3").Value = -519.25
$ws.Range("N3").Value = -51134.5
$ws.Range("H20").Value = 1228.2106
$ws.Range("I20").Value = 937.41174
$ws.Range("J20").Value = 3700
$ws.Range("K20").Value = 937.41174
$ws.Range("L20").Value = 3700
$ws.Range("M20").Value = -690.41174
$ws.Range("N20").Value = -4194
$ws.Range("H94").Value = 41667116
$ws.Range("I94").Value = 41667116
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 41667116
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -41666665
$ws.Range("N94").ClearContents()
$ws.Range("H99").Value = 50001360
$ws.Range("I99").Value = 55556844
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 55556844
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -55555346
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 125003520
$ws.Range("I105").Value = 125003520
$ws.Range("K105").Value = 125003520
$ws.Range("M105").Value = -125001773
$ws.Range("H107").Value = 1154.7
$ws.Range("I107").Value = 1050.0625
$ws.Range("J107").Value = 1573.25
$ws.Range("K107").Value = 1050.0625
$ws.Range("L107").Value = 1573.25
$ws.Range("M107").Value = 869.9375
$ws.Range("N107").Value = -5413.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1237.2167
$ws.Range("I31").Value = 1176.9636
$ws.Range("K31").Value = 1176.9636
$ws.Range("M31").Value = -881.9636
$ws.Range("H34").Value = 1237.2167
$ws.Range("I34").Value = 1176.9636
$ws.Range("K34").Value = 1176.9636
$ws.Range("M34").Value = -974.9636
$ws.Range("H99").Value = 2118.6667
$ws.Range("J99").Value = 2350
$ws.Range("L99").Value = 2350
$ws.Range("N99").Value = -5346
$ws.Range("H126").Value = 2118.6667
$ws.Range("J126").Value = 2350
$ws.Range("L126").Value = 7050
$ws.Range("N126").Value = -11990
$ws.Range("H132").Value = 5299.968
$ws.Range("I132").Value = 6185.3184
$ws.Range("J132").Value = 3135.7778
$ws.Range("K132").Value = 18555.9552
$ws.Range("L132").Value = 9407.3334
$ws.Range("M132").Value = -16025.9552
$ws.Range("N132").Value = -14467.3334
$ws.Range("H134").Value = 2635.2666
$ws.Range("I134").Value = 3259.4443
$ws.Range("J134").Value = 1699
$ws.Range("K134").Value = 9778.332900000001
$ws.Range("L134").Value = 5097
$ws.Range("M134").Value = -7243.332900000001
$ws.Range("N134").Value = -10167
$ws.Range("H141").Value = 27518.666
$ws.Range("J141").Value = 28622.4
$ws.Range("L141").Value = 28622.4
$ws.Range("N141").Value = -38982.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12660464
$ws.Range("J131").Value = 2335.0266
$ws.Range("L131").Value = 7005.0798
$ws.Range("N131").Value = -17085.0798

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 40912784
$ws.Range("I70").Value = 27781756
$ws.Range("J70").Value = 100002400
$ws.Range("K70").Value = 27781756
$ws.Range("L70").Value = 100002400
$ws.Range("M70").Value = -27781486
$ws.Range("N70").Value = -100002940
$ws.Range("H73").Value = 40912784
$ws.Range("I73").Value = 27781756
$ws.Range("J73").Value = 100002400
$ws.Range("K73").Value = 27781756
$ws.Range("L73").Value = 100002400
$ws.Range("M73").Value = -27780820
$ws.Range("N73").Value = -100004272
$ws.Range("H80").Value = 4750.7144
$ws.Range("I80").Value = 3163.75
$ws.Range("K80").Value = 3163.75
$ws.Range("M80").Value = -2165.75
$ws.Range("H83").Value = 4750.7144
$ws.Range("I83").Value = 3163.75
$ws.Range("K83").Value = 15818.75
$ws.Range("M83").Value = -10826.75
$ws.Range("H132").Value = 2407.889
$ws.Range("I132").Value = 2018.1904
$ws.Range("J132").Value = 2953.4666
$ws.Range("K132").Value = 6054.5712
$ws.Range("L132").Value = 8860.399800000001
$ws.Range("M132").Value = -3524.5712
$ws.Range("N132").Value = -13920.3998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H93").Value = 1200
$ws.Range("J93").Value = 1200
$ws.Range("L93").Value = 1200
$ws.Range("N93").Value = -3696
$ws.Range("H122").Value = 15633829
$ws.Range("I122").Value = 27791164
$ws.Range("J122").Value = 2968.5715
$ws.Range("K122").Value = 83373492
$ws.Range("L122").Value = 8905.7145
$ws.Range("M122").Value = -83371042
$ws.Range("N122").Value = -13805.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 5000
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 5000
$ws.Range("N69").Value = -6498
$ws.Range("H72").Value = 5000
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 15000
$ws.Range("N72").Value = -22488
$ws.Range("H96").Value = 2911.2307
$ws.Range("I96").Value = 2899.6365
$ws.Range("J96").Value = 2975
$ws.Range("K96").Value = 2899.6365
$ws.Range("L96").Value = 2975
$ws.Range("M96").Value = -1526.6365
$ws.Range("N96").Value = -5721
$ws.Range("H132").Value = 2810.2
$ws.Range("J132").Value = 2581.0908
$ws.Range("L132").Value = 7743.2724
$ws.Range("N132").Value = -12803.2724
$ws.Range("H136").Value = 627
$ws.Range("I136").Value = 427.73685
$ws.Range("K136").Value = 1283.21055
$ws.Range("M136").Value = 1266.78945
